$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 11 new rows before the old row 29 (the last data row / thick-bottom row),
#    pushing it (and everything below) down to row 40. Excel auto-extends the
#    SUM(...) formulas and the thick-bottom formatting of that row because the
#    insertion happens inside the referenced ranges.
$ws.Rows("29:39").Insert(-4121, 0)

# 2. Re-apply the correct (pre-existing) cell formatting to the newly inserted rows.
#    New rows 29-31 should pick up the formatting used by the earlier "first group"
#    of blank data rows (e.g. row 17), and new rows 32-39 should pick up the
#    formatting used by the later "second group" of blank data rows (e.g. row 28,
#    now shifted to row 39 -... use an untouched donor row instead, row 17/28 area
#    is unaffected by the insert so it is safe to use directly).
$ws.Range("B17:K17").Copy()
$ws.Range("B29:K31").PasteSpecial(-4122)

$ws.Range("B28:K28").Copy()
$ws.Range("B32:K39").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3. Set the row heights to match the donor rows (22.8) for the newly created rows.
$ws.Rows("29:39").RowHeight = 22.8

# 4. Fill in column A (row index numbers) for the newly inserted rows: 25-35,
#    and renumber the old last-data row (shifted from row 29 to row 40) to 36.
$ws.Range("A29").Value = 25
$ws.Range("A30").Value = 26
$ws.Range("A31").Value = 27
$ws.Range("A32").Value = 28
$ws.Range("A33").Value = 29
$ws.Range("A34").Value = 30
$ws.Range("A35").Value = 31
$ws.Range("A36").Value = 32
$ws.Range("A37").Value = 33
$ws.Range("A38").Value = 34
$ws.Range("A39").Value = 35
$ws.Range("A40").Value = 36

# 5. Update the print area defined name to cover the new extent of the sheet.
$ws.PageSetup.PrintArea = "`$B`$1:`$K`$44"

# 6. Update the active selection to reflect where the author left off editing.
$ws.Range("F41").Select()

$wb.Save()
